# Updated symbol list on Wed Jan 25 17:23:58 UTC 2023 with GitHub Actions
# Refreshes the Price / Volume(1h) columns for every coin row, and
# re-syncs rows 8-23 where the ranking shuffled (each row's Coin/Link
# pair moved up/down relative to its neighbours).
#
# Note: Price/Volume cells hold numeric-looking text (e.g. "301.79",
# "-3.23%"), not real numbers - a leading "'" forces Excel to store the
# literal string instead of auto-converting it to a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.79"
$ws.Range("E2").Value = "'-3.23%"
$ws.Range("D3").Value = "'35.34"
$ws.Range("E3").Value = "'-0.15%"
$ws.Range("D4").Value = "'5.061"
$ws.Range("E4").Value = "'-0.42%"
$ws.Range("D5").Value = "'0.07928"
$ws.Range("E5").Value = "'-2.94%"
$ws.Range("D6").Value = "'1.895"
$ws.Range("E6").Value = "'-9.61%"
$ws.Range("D7").Value = "'7.773"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9275"
$ws.Range("E8").Value = "'0.18%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1380"
$ws.Range("E9").Value = "'31.28%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1905"
$ws.Range("E10").Value = "'-1.03%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09156"
$ws.Range("E11").Value = "'-0.62%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03432"
$ws.Range("E12").Value = "'-5.46%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09922"
$ws.Range("E13").Value = "'0.11%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001389"
$ws.Range("E14").Value = "'-3.08%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04512"
$ws.Range("E15").Value = "'-0.62%"
$ws.Range("D16").Value = "'0.005831"
$ws.Range("E16").Value = "'1.13%"
$ws.Range("D17").Value = "'3.530"
$ws.Range("E17").Value = "'1.53%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.045"
$ws.Range("E18").Value = "'-2.10%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.926"
$ws.Range("E19").Value = "'4.02%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3407"
$ws.Range("E20").Value = "'-0.15%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1294"
$ws.Range("E21").Value = "'-0.59%"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "'5.043"
$ws.Range("E22").Value = "'-0.96%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.2398"
$ws.Range("E23").Value = "'8.26%"
$ws.Range("D24").Value = "'0.001212"
$ws.Range("E24").Value = "'-1.27%"
$ws.Range("E25").Value = "'-0.23%"
$ws.Range("E26").Value = "'-1.57%"
$ws.Range("D27").Value = "'0.0003001"
$ws.Range("E27").Value = "'-32.53%"
$ws.Range("D39").Value = "'0.01864"
$ws.Range("E39").Value = "'-5.35%"
$ws.Range("D40").Value = "'0.04774"
$ws.Range("E40").Value = "'-2.57%"
$ws.Range("D41").Value = "'0.007332"
$ws.Range("E41").Value = "'-3.19%"
$ws.Range("D42").Value = "'0.009629"
$ws.Range("E42").Value = "'9.73%"
$ws.Range("E43").Value = "'-4.09%"
$ws.Range("E44").Value = "'-2.41%"
$ws.Range("D45").Value = "'0.01100"
$ws.Range("E45").Value = "'-5.66%"
$ws.Range("D46").Value = "'0.00006246"
$ws.Range("E46").Value = "'-5.69%"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("D48").Value = "'64.66"
$ws.Range("E48").Value = "'-64.92%"
$ws.Range("E49").Value = "'10.63%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.03%"
